$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Row 69
$ws.Cells.Item(69, 1).Value = 2
$ws.Cells.Item(69, 2).Value = "Five Families of Crime"
$ws.Cells.Item(69, 3).Value = "Poison Thanos"
$ws.Cells.Item(69, 4).Value = "Poisons|Sinister Six|Mojoverse|Skrulls"
$ws.Cells.Item(69, 5).Value = "Hellfire Cult"
$ws.Cells.Item(69, 6).Value = "Soulsword Colossus (SW2)|Skaar, Son of Hulk (WW)|Sabretooth (V)|Gladiator Hulk (WW)|Hercules (CW)"
$ws.Cells.Item(69, 7).Value = 1
$ws.Cells.Item(69, 8).Value = "50|44"
$ws.Cells.Item(69, 9).Value = "not really"
$ws.Cells.Item(69, 11).Value = "6 escapees, lots of rare heroes acquired. Skaar and Sabretooth are a good combo."

# Row 70
$ws.Cells.Item(70, 1).Value = 2
$ws.Cells.Item(70, 2).Value = "Earthquake Drains the Ocean"
$ws.Cells.Item(70, 3).Value = "Mojo"
$ws.Cells.Item(70, 4).Value = "Mojoverse|The Deadlands|Zola's Creations"
$ws.Cells.Item(70, 5).Value = "Circus of Crime"
$ws.Cells.Item(70, 6).Value = "Warlock (NM)|Spider-Man (B)|Medusa (ROK)|Greithoth, Breaker of Wills (FI)|Captain Marvel, Agent of S.H.I.E.L.D. (R)"
$ws.Cells.Item(70, 7).Value = 1
$ws.Cells.Item(70, 11).Value = "Five villains escaped early before stabilizing. Spiral ruined two whole turns for both players. Managed to dodge the living dead rising again and spiderman did quite some work."
$ws.Cells.Item(70, 8).Value = "58|55"
$ws.Cells.Item(70, 9).Value = "yes"

# Row 71
$ws.Cells.Item(71, 1).Value = 2
$ws.Cells.Item(71, 2).Value = "Deadpool Wants a Chimichanga"
$ws.Cells.Item(71, 3).Value = "Supreme Intelligence of the Kree"
$ws.Cells.Item(71, 4).Value = "Kree Starforce|Utopolis"
$ws.Cells.Item(71, 5).Value = "Thor Corps"
$ws.Cells.Item(71, 6).Value = "Venompool (VE)|Spider-Woman (PTT)|Ultimate Spider-Man (SW1)|Rick Jones (WW)|Namor, the Sub-Mariner (SW1)"
$ws.Cells.Item(71, 7).Value = 0
$ws.Cells.Item(71, 11).Value = "One turn short to get the last tactic before the villain deck ran out. Thor Corps was a major blocker."
$ws.Cells.Item(71, 8).Value = "32|38"
$ws.Cells.Item(71, 9).Value = "yes"

$ws.Range("H72").Select()
